$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (Förändrad) for rows 2-5 from date serial 45184 to 45185
# (2023-09-15 -> 2023-09-16), keeping existing cell formatting/style intact.
$ws.Range("C2:C5").Value = 45185
